$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 7.739652666666667
$ws.Range("H2").Value = 23.218958
$ws.Range("I2").Value = 0.6488398532974882
$ws.Range("J2").Value = 0.6488398532974882
$ws.Range("M2").Value = 8.489835333333334
$ws.Range("N2").Value = 25.469506
$ws.Range("O2").Value = 0.2075776945087381
$ws.Range("P2").Value = 0.2075776945087381
$ws.Range("Q2").Value = 65.70837667719422
$ws.Range("R2").Value = 591.375390094748
$ws.Range("S2").Value = 0.1346846808528805
$ws.Range("T2").Value = 0.1346846808528805
$ws.Range("G3").Value = 7.739652666666667
$ws.Range("H3").Value = 23.218958
$ws.Range("I3").Value = 0.6488398532974882
$ws.Range("J3").Value = 0.6488398532974882
$ws.Range("O3").Value = 0.3214784855238645
$ws.Range("P3").Value = 0.3214784855238645
$ws.Range("Q3").Value = 101.7634841277553
$ws.Range("R3").Value = 915.871357149798
$ws.Range("S3").Value = 0.2085880533856029
$ws.Range("T3").Value = 0.2085880533856029
$ws.Range("G4").Value = 7.739652666666667
$ws.Range("H4").Value = 23.218958
$ws.Range("I4").Value = 0.6488398532974882
$ws.Range("J4").Value = 0.6488398532974882
$ws.Range("M4").Value = 5.630791333333334
$ws.Range("N4").Value = 16.892374
$ws.Range("O4").Value = 0.1376736576555254
$ws.Range("P4").Value = 0.1376736576555254
$ws.Range("Q4").Value = 43.58036915847689
$ws.Range("R4").Value = 392.223322426292
$ws.Range("S4").Value = 0.08932815583613973
$ws.Range("T4").Value = 0.08932815583613972
$ws.Range("G5").Value = 7.739652666666667
$ws.Range("H5").Value = 23.218958
$ws.Range("I5").Value = 0.6488398532974882
$ws.Range("J5").Value = 0.6488398532974882
$ws.Range("M5").Value = 6.738585333333333
$ws.Range("N5").Value = 20.215756
$ws.Range("O5").Value = 0.1647593802263456
$ws.Range("P5").Value = 0.1647593802263456
$ws.Range("Q5").Value = 52.15430994469423
$ws.Range("R5").Value = 469.388789502248
$ws.Range("S5").Value = 0.1069024520954471
$ws.Range("T5").Value = 0.1069024520954471
$ws.Range("G6").Value = 7.739652666666667
$ws.Range("H6").Value = 23.218958
$ws.Range("I6").Value = 0.6488398532974882
$ws.Range("J6").Value = 0.6488398532974882
$ws.Range("M6").Value = 0.9376886666666667
$ws.Range("N6").Value = 2.813066
$ws.Range("O6").Value = 0.02292662271427322
$ws.Range("P6").Value = 0.02292662271427321
$ws.Range("Q6").Value = 7.257384589469779
$ws.Range("R6").Value = 65.316461305228
$ws.Range("S6").Value = 0.01487570651853589
$ws.Range("T6").Value = 0.01487570651853589
$ws.Range("G7").Value = 7.739652666666667
$ws.Range("H7").Value = 23.218958
$ws.Range("I7").Value = 0.6488398532974882
$ws.Range("J7").Value = 0.6488398532974882
$ws.Range("M7").Value = 5.954327333333333
$ws.Range("N7").Value = 17.862982
$ws.Range("O7").Value = 0.1455841593712531
$ws.Range("P7").Value = 0.1455841593712531
$ws.Range("Q7").Value = 46.08442542363955
$ws.Range("R7").Value = 414.759828812756
$ws.Range("S7").Value = 0.09446080460888202
$ws.Range("T7").Value = 0.09446080460888202
$ws.Range("I8").Value = 0.3053032463428815
$ws.Range("J8").Value = 0.3053032463428815
$ws.Range("M8").Value = 8.489835333333334
$ws.Range("N8").Value = 25.469506
$ws.Range("O8").Value = 0.2075776945087381
$ws.Range("P8").Value = 0.2075776945087381
$ws.Range("Q8").Value = 30.91823137792134
$ws.Range("R8").Value = 278.264082401292
$ws.Range("S8").Value = 0.06337414400188868
$ws.Range("T8").Value = 0.06337414400188866
$ws.Range("I9").Value = 0.3053032463428815
$ws.Range("J9").Value = 0.3053032463428815
$ws.Range("O9").Value = 0.3214784855238645
$ws.Range("P9").Value = 0.3214784855238645
$ws.Range("S9").Value = 0.09814842525982885
$ws.Range("T9").Value = 0.09814842525982885
$ws.Range("I10").Value = 0.3053032463428815
$ws.Range("J10").Value = 0.3053032463428815
$ws.Range("M10").Value = 5.630791333333334
$ws.Range("N10").Value = 16.892374
$ws.Range("O10").Value = 0.1376736576555254
$ws.Range("P10").Value = 0.1376736576555254
$ws.Range("Q10").Value = 20.50618209298534
$ws.Range("R10").Value = 184.555638836868
$ws.Range("S10").Value = 0.04203221461813041
$ws.Range("T10").Value = 0.04203221461813041
$ws.Range("I11").Value = 0.3053032463428815
$ws.Range("J11").Value = 0.3053032463428815
$ws.Range("M11").Value = 6.738585333333333
$ws.Range("N11").Value = 20.215756
$ws.Range("O11").Value = 0.1647593802263456
$ws.Range("P11").Value = 0.1647593802263456
$ws.Range("Q11").Value = 24.54053963542134
$ws.Range("R11").Value = 220.864856718792
$ws.Range("S11").Value = 0.05030157364854446
$ws.Range("T11").Value = 0.05030157364854446
$ws.Range("I12").Value = 0.3053032463428815
$ws.Range("J12").Value = 0.3053032463428815
$ws.Range("M12").Value = 0.9376886666666667
$ws.Range("N12").Value = 2.813066
$ws.Range("O12").Value = 0.02292662271427322
$ws.Range("P12").Value = 0.02292662271427321
$ws.Range("Q12").Value = 3.414868960134667
$ws.Range("R12").Value = 30.733820641212
$ws.Range("S12").Value = 0.006999572342346058
$ws.Range("T12").Value = 0.006999572342346057
$ws.Range("I13").Value = 0.3053032463428815
$ws.Range("J13").Value = 0.3053032463428815
$ws.Range("M13").Value = 5.954327333333333
$ws.Range("N13").Value = 17.862982
$ws.Range("O13").Value = 0.1455841593712531
$ws.Range("P13").Value = 0.1455841593712531
$ws.Range("Q13").Value = 21.68443355656933
$ws.Range("R13").Value = 195.159902009124
$ws.Range("S13").Value = 0.04444731647214301
$ws.Range("T13").Value = 0.04444731647214301
$ws.Range("G14").Value = 0.5470016666666667
$ws.Range("H14").Value = 1.641005
$ws.Range("I14").Value = 0.04585690035963046
$ws.Range("J14").Value = 0.04585690035963046
$ws.Range("M14").Value = 8.489835333333334
$ws.Range("N14").Value = 25.469506
$ws.Range("O14").Value = 0.2075776945087381
$ws.Range("P14").Value = 0.2075776945087381
$ws.Range("Q14").Value = 4.64395407705889
$ws.Range("R14").Value = 41.79558669353
$ws.Range("S14").Value = 0.009518869653969017
$ws.Range("T14").Value = 0.009518869653969014
$ws.Range("G15").Value = 0.5470016666666667
$ws.Range("H15").Value = 1.641005
$ws.Range("I15").Value = 0.04585690035963046
$ws.Range("J15").Value = 0.04585690035963046
$ws.Range("O15").Value = 0.3214784855238645
$ws.Range("P15").Value = 0.3214784855238645
$ws.Range("Q15").Value = 7.192156782878334
$ws.Range("R15").Value = 64.729411045905
$ws.Range("S15").Value = 0.01474200687843276
$ws.Range("T15").Value = 0.01474200687843276
$ws.Range("G16").Value = 0.5470016666666667
$ws.Range("H16").Value = 1.641005
$ws.Range("I16").Value = 0.04585690035963046
$ws.Range("J16").Value = 0.04585690035963046
$ws.Range("M16").Value = 5.630791333333334
$ws.Range("N16").Value = 16.892374
$ws.Range("O16").Value = 0.1376736576555254
$ws.Range("P16").Value = 0.1376736576555254
$ws.Range("Q16").Value = 3.080052243985556
$ws.Range("R16").Value = 27.72047019587
$ws.Range("S16").Value = 0.006313287201255306
$ws.Range("T16").Value = 0.006313287201255304
$ws.Range("G17").Value = 0.5470016666666667
$ws.Range("H17").Value = 1.641005
$ws.Range("I17").Value = 0.04585690035963046
$ws.Range("J17").Value = 0.04585690035963046
$ws.Range("M17").Value = 6.738585333333333
$ws.Range("N17").Value = 20.215756
$ws.Range("O17").Value = 0.1647593802263456
$ws.Range("P17").Value = 0.1647593802263456
$ws.Range("Q17").Value = 3.686017408308889
$ws.Range("R17").Value = 33.17415667478
$ws.Range("S17").Value = 0.007555354482353998
$ws.Range("T17").Value = 0.007555354482353997
$ws.Range("G18").Value = 0.5470016666666667
$ws.Range("H18").Value = 1.641005
$ws.Range("I18").Value = 0.04585690035963046
$ws.Range("J18").Value = 0.04585690035963046
$ws.Range("M18").Value = 0.9376886666666667
$ws.Range("N18").Value = 2.813066
$ws.Range("O18").Value = 0.02292662271427322
$ws.Range("P18").Value = 0.02292662271427321
$ws.Range("Q18").Value = 0.5129172634811112
$ws.Range("R18").Value = 4.61625537133
$ws.Range("S18").Value = 0.001051343853391268
$ws.Range("T18").Value = 0.001051343853391267
$ws.Range("G19").Value = 0.5470016666666667
$ws.Range("H19").Value = 1.641005
$ws.Range("I19").Value = 0.04585690035963046
$ws.Range("J19").Value = 0.04585690035963046
$ws.Range("M19").Value = 5.954327333333333
$ws.Range("N19").Value = 17.862982
$ws.Range("O19").Value = 0.1455841593712531
$ws.Range("P19").Value = 0.1455841593712531
$ws.Range("Q19").Value = 3.257026975212222
$ws.Range("R19").Value = 29.31324277691
$ws.Range("S19").Value = 0.006676038290228116
$ws.Range("T19").Value = 0.006676038290228115
